$wb = $excel.ActiveWorkbook

# --- Sheet "Results" (sheet1) ---
$wsResults = $wb.Worksheets.Item("Results")
$wsResults.Range("C2").Value = 3501
$wsResults.Range("D3").Value = 3794
$wsResults.Range("D4").Value = 3612
$wsResults.Range("C5").Value = 1412
$wsResults.Range("D5").Value = 3557
$wsResults.Range("D6").Value = 10571
$wsResults.Range("C7").Value = 12482
$wsResults.Range("D7").Value = 6944
$wsResults.Range("C8").Value = 12931
$wsResults.Range("D9").Value = 4991
$wsResults.Range("E9").Value = 1882
$wsResults.Range("C10").Value = 7711
$wsResults.Range("C11").Value = 18401
$wsResults.Range("D11").Value = 20878
$wsResults.Range("C12").Value = 17755

# --- Sheet "PRIOR" (sheet2) ---
$wsPrior = $wb.Worksheets.Item("PRIOR")
$wsPrior.Range("B2").Value = 58709
$wsPrior.Range("C2").Value = 35755
$wsPrior.Range("D2").Value = 27462

# --- Sheet "LEAD" (sheet3) ---
$wsLead = $wb.Worksheets.Item("LEAD")
$wsLead.Range("B2").Value = 7782
$wsLead.Range("C2").Value = 5205
$wsLead.Range("B6").Value = 8072
$wsLead.Range("C6").Value = 9130
$wsLead.Range("D6").Value = 8606
$wsLead.Range("C7").Value = 278
$wsLead.Range("D12").Value = 236
$wsLead.Range("D13").Value = 0
$wsLead.Range("B16").Value = 13153
$wsLead.Range("C16").Value = 16096
$wsLead.Range("B17").Value = 48571
$wsLead.Range("C17").Value = 43428
$wsLead.Range("D17").Value = 15606
$wsLead.Range("C19").Value = 811
$wsLead.Range("D20").Value = 4319

# Column A width change on LEAD sheet (target stored width is exactly 25;
# Excel's ColumnWidth<->stored-width conversion rounds to the nearest pixel,
# so 24.14 characters maps to a stored width of 25.0 in the XML)
$wsLead.Columns.Item(1).ColumnWidth = 24.14
